$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    if ($value -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

# Row 2
Set-TextCell 2 4 "43.444.83"
Set-TextCell 2 5 "  -1.16%  "

# Row 3
Set-TextCell 3 4 "2.375.14"
Set-TextCell 3 5 "  +5.79%  "

# Row 4
Set-TextCell 4 5 "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 2).Value = "XRP"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell 5 4 "0.657"
Set-TextCell 5 5 "  +2.64%  "

# Row 6
$ws.Cells.Item(6, 2).Value = "BNB"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell 6 4 "233.88"
Set-TextCell 6 5 "  +1.12%  "

# Row 7
Set-TextCell 7 4 "69.73"
Set-TextCell 7 5 "  +11.01%  "

# Row 8
Set-TextCell 8 5 "  +0.13%  "

# Row 9
Set-TextCell 9 4 "0.461"
Set-TextCell 9 5 "  +3.62%  "

# Row 10
Set-TextCell 10 4 "0.0968"
Set-TextCell 10 5 "  -0.92%  "

# Row 11
Set-TextCell 11 4 "57.21"
Set-TextCell 11 5 "  +0.24%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell 12 4 "2.735.73"
Set-TextCell 12 5 "  +5.95%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "Avalanche"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 13 4 "26.20"
Set-TextCell 13 5 "  +0.39%  "

# Row 14
Set-TextCell 14 5 "  +0.02%  "

# Row 15
Set-TextCell 15 4 "15.67"

# Row 16
Set-TextCell 16 4 "6.23"
Set-TextCell 16 5 "  +2.42%  "

# Row 17
Set-TextCell 17 4 "0.851"
Set-TextCell 17 5 "  +3.24%  "

# Row 18
Set-TextCell 18 4 "2.388.00"
Set-TextCell 18 5 "  +5.90%  "

# Row 19
Set-TextCell 19 4 "43.498.96"
Set-TextCell 19 5 "  -0.81%  "

# Row 20
Set-TextCell 20 5 "  +0.85%  "

# Row 21
Set-TextCell 21 4 "6.33"
Set-TextCell 21 5 "  +4.85%  "

# Row 22
Set-TextCell 22 4 "74.10"
Set-TextCell 22 5 "  +2.19%  "

# Row 23
Set-TextCell 23 4 "251.07"
Set-TextCell 23 5 "  +1.47%  "

# Row 24
Set-TextCell 24 5 "  +18.31%  "

# Row 25
Set-TextCell 25 5 "  +0.03%  "

# Row 26
Set-TextCell 26 4 "2.45"
Set-TextCell 26 5 "  +1.93%  "

# Row 27
Set-TextCell 27 5 "  +2.37%  "

# Row 28
Set-TextCell 28 4 "22.89"
Set-TextCell 28 5 "  +9.30%  "

# Row 29
Set-TextCell 29 5 "  +2.24%  "

# Row 30
Set-TextCell 30 4 "174.21"
Set-TextCell 30 5 "  +1.32%  "

# Row 31
Set-TextCell 31 4 "1.54"
Set-TextCell 31 5 "  +9.91%  "

# Row 32
Set-TextCell 32 5 "  -8.63%  "

# Row 33
Set-TextCell 33 4 "0.128"
Set-TextCell 33 5 "  +2.43%  "

# Row 34
Set-TextCell 34 4 "4.97"
Set-TextCell 34 5 "  +4.40%  "

# Row 35
Set-TextCell 35 5 "  +1.12%  "

# Row 36
Set-TextCell 36 4 "5.08"
Set-TextCell 36 5 "  +3.26%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "THORChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell 37 4 "6.59"
Set-TextCell 37 5 "  +3.18%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "LidoDAOToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 38 4 "2.44"
Set-TextCell 38 5 "  +7.70%  "

# Row 39
Set-TextCell 39 4 "3.63"
Set-TextCell 39 5 "  +0.05%  "

# Row 40
Set-TextCell 40 5 "  +1.17%  "

# Row 41
Set-TextCell 41 5 "  +4.70%  "

# Row 42
Set-TextCell 42 5 "  -0.18%  "

# Row 43
Set-TextCell 43 4 "18.54"
Set-TextCell 43 5 "  +9.26%  "

# Row 44
Set-TextCell 44 4 "1.18"
Set-TextCell 44 5 "  +10.68%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "TrustWalletToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 45 4 "1.22"
Set-TextCell 45 5 "  +2.38%  "

# Row 46
Set-TextCell 46 5 "  +4.55%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 47 4 "99.20"
Set-TextCell 47 5 "  +2.21%  "

# Row 48
Set-TextCell 48 4 "0.0947"
Set-TextCell 48 5 "  +0.41%  "

# Row 49
Set-TextCell 49 4 "1.450.69"
Set-TextCell 49 5 "  +0.91%  "

# Row 50
Set-TextCell 50 4 "2.605.88"
Set-TextCell 50 5 "  +6.17%  "

# Row 51
Set-TextCell 51 5 "  -0.66%  "
